$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New day of code-coverage numbers (column D = 2016-04-24) for each tracked
# category, entered below the existing column B figures.
$ws.Range("D2").Value = 0.338    # entity
$ws.Range("D3").Value = 0        # controller
$ws.Range("D4").Value = 0        # GUI
$ws.Range("D5").Value = 0.517    # move
$ws.Range("D6").Value = 0        # view
$ws.Range("D8").Value = 0.171    # total

# Match the percentage number format already used in column B (numFmtId 10,
# "0.00%") so the new cells share the same style instead of minting a new one.
$ws.Range("D2:D6").NumberFormat = "0.00%"
$ws.Range("D8").NumberFormat = "0.00%"

# Leave the selection where the author last clicked while entering this data.
$ws.Range("E5").Select()
